$wb = $excel.ActiveWorkbook

# Fix typo "witnessned" -> "witnessed" on Ive_Second_Phrase sheet (A5)
$ws6 = $wb.Worksheets.Item("Ive_Second_Phrase")
$ws6.Range("A5").Value = "witnessed"

# Clear the "Punctuation" NA values in rows 3-8 of Your_Second_Phrase
$ws4 = $wb.Worksheets.Item("Your_Second_Phrase")
$ws4.Range("D3:D8").ClearContents()

# Clear the "Punctuation" NA values in rows 3-7 of Youre_Second_Phrase
$ws5 = $wb.Worksheets.Item("Youre_Second_Phrase")
$ws5.Range("D3:D7").ClearContents()

# Update selections / active sheet to match final workbook state
$ws6.Range("A6").Select()
$ws5.Range("D7").Select()
$ws4.Range("D3:D8").Select()

$ws4.Activate()
